$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are plain decimal-looking strings (e.g. '1.001',
# '121.50', '5.500', '0.00001057'). Column D stores these as literal text
# (matching coinranking.com's raw formatting), so force each such cell to
# Text format before assigning the string -- otherwise Excel would silently
# reinterpret them as numbers and lose the exact formatting (trailing zeros,
# scientific notation, etc.).
$textCells = @(
    'D4',
    'D5',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D13',
    'D14',
    'D15',
    'D17',
    'D18',
    'D19',
    'D21',
    'D22',
    'D24',
    'D25',
    'D26',
    'D27',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.811.38'
$ws.Range('E2').Value = '  -3.04%  '
$ws.Range('D3').Value = '1.791.94'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '315.87'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D7').Value = '0.5351'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = '0.3823'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').Value = '0.07438'
$ws.Range('D10').Value = '41.44'
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('D11').Value = '1.085'
$ws.Range('E11').Value = '  -2.80%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '6.194'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = '7.427'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').Value = '20.31'
$ws.Range('E15').Value = '  -2.38%  '
$ws.Range('D16').Value = '1.795.35'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').Value = '88.26'
$ws.Range('E17').Value = '  -2.40%  '
$ws.Range('D18').Value = '0.00001057'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').Value = '0.06513'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '17.28'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = '5.959'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').Value = '27.858.70'
$ws.Range('D24').Value = '11.12'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = '2.093'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = '157.21'
$ws.Range('E26').Value = '  -1.68%  '
$ws.Range('D27').Value = '20.17'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '1.994.55'
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('D29').Value = '2.322'
$ws.Range('E29').Value = '  -2.47%  '
$ws.Range('D30').Value = '121.50'
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.1093'
$ws.Range('E31').Value = '  +3.84%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '1.104'
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').Value = '3.649'
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('D34').Value = '5.500'
$ws.Range('E34').Value = '  -2.71%  '
$ws.Range('D35').Value = '0.06941'
$ws.Range('E35').Value = '  +7.60%  '
$ws.Range('D36').Value = '0.2186'
$ws.Range('E36').Value = '  -2.87%  '
$ws.Range('D37').Value = '0.02269'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('D38').Value = '5.037'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '11.37'
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').Value = '8.407'
$ws.Range('E40').Value = '  -5.66%  '
$ws.Range('D41').Value = '0.6092'
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('D42').Value = '1.167'
$ws.Range('E42').Value = '  -4.99%  '
$ws.Range('D43').Value = '1.409'
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('D44').Value = '13.28'
$ws.Range('E44').Value = '  -0.29%  '
$ws.Range('D45').Value = '3.677'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').Value = '0.5687'
$ws.Range('E46').Value = '  -3.24%  '
$ws.Range('D47').Value = '125.04'
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('D48').Value = '1.907'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').Value = '1.170'
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('D50').Value = '0.06801'
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('D51').Value = '71.24'
$ws.Range('E51').Value = '  -1.90%  '
